$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F5").Value = -4
$ws.Range("F7").Value = -4
$ws.Range("F8").Value = 2
$ws.Range("F9").Value = -4
$ws.Range("F10").Value = -5
$ws.Range("F15").Value = -4
$ws.Range("F16").Value = -3
$ws.Range("F17").Value = -10
$ws.Range("F18").Value = -7
$ws.Range("F22").Value = -6
$ws.Range("F24").Value = -5
$ws.Range("F25").Value = -5
$ws.Range("F26").Value = -8
$ws.Range("F30").Value = -1
$ws.Range("F31").Value = 4
$ws.Range("F34").Value = -3
$ws.Range("F35").Value = 0
$ws.Range("F39").Value = -2
$ws.Range("F41").Value = -6
$ws.Range("F43").Value = 0
$ws.Range("F47").Value = -5
$ws.Range("F49").Value = 3
$ws.Range("F54").Value = 4
$ws.Range("F66").Value = 2
$ws.Range("F67").Value = 0
$ws.Range("F70").Value = 1
$ws.Range("F71").Value = -5
$ws.Range("F72").Value = -2
